$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $savedStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $savedStyle
}

$ws.Range('D2').Value = '37.727.81'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').Value = '2.073.39'
$ws.Range('E3').Value = '  -1.30%  '
$ws.Range('E4').Value = '  -0.02%  '
Set-TextValue 'D5' '232.71'
$ws.Range('E5').Value = '  -0.57%  '
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('E7').Value = '  -0.05%  '
Set-TextValue 'D8' '58.32'
$ws.Range('E8').Value = '  +1.26%  '
$ws.Range('E9').Value = '  +1.12%  '
$ws.Range('E10').Value = '  +1.00%  '
$ws.Range('E11').Value = '  +3.38%  '
$ws.Range('D12').Value = '2.379.96'
$ws.Range('E12').Value = '  -1.39%  '
Set-TextValue 'D13' '14.70'
$ws.Range('E13').Value = '  +2.18%  '
Set-TextValue 'D14' '20.83'
$ws.Range('E14').Value = '  -1.56%  '
$ws.Range('E15').Value = '  -0.89%  '
$ws.Range('E16').Value = '  +2.73%  '
$ws.Range('D17').Value = '2.107.16'
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('D18').Value = '37.667.56'
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('E19').Value = '  -0.53%  '
Set-TextValue 'D20' '70.98'
$ws.Range('E20').Value = '  +1.24%  '
$ws.Range('D21').Value = '0.0₃0831'
$ws.Range('E21').Value = '  +1.50%  '
Set-TextValue 'D22' '228.04'
$ws.Range('E22').Value = '  +0.60%  '
$ws.Range('E23').Value = '  +0.00%  '
Set-TextValue 'D24' '2.37'
$ws.Range('E24').Value = '  -2.13%  '
$ws.Range('E25').Value = '  -0.20%  '
Set-TextValue 'D26' '170.68'
$ws.Range('E26').Value = '  +0.98%  '
$ws.Range('E27').Value = '  +4.46%  '
Set-TextValue 'D28' '9.00'
$ws.Range('E28').Value = '  +1.10%  '
$ws.Range('E29').Value = '  +0.43%  '
$ws.Range('E30').Value = '  -1.43%  '
$ws.Range('E31').Value = '  +2.67%  '
$ws.Range('E32').Value = '  +1.56%  '
$ws.Range('E33').Value = '  +1.49%  '
Set-TextValue 'D34' '4.64'
$ws.Range('E34').Value = '  +2.32%  '
$ws.Range('E35').Value = '  -2.96%  '
$ws.Range('E36').Value = '  +0.30%  '
Set-TextValue 'D37' '3.38'
$ws.Range('E37').Value = '  -1.15%  '
$ws.Range('E38').Value = '  -0.08%  '
Set-TextValue 'D39' '5.30'
$ws.Range('E39').Value = '  -2.07%  '
Set-TextValue 'D40' '100.76'
$ws.Range('E40').Value = '  +4.69%  '
Set-TextValue 'D41' '0.0973'
$ws.Range('E41').Value = '  -3.60%  '
$ws.Range('E42').Value = '  -1.89%  '
Set-TextValue 'D43' '0.0214'
$ws.Range('E43').Value = '  +0.94%  '
$ws.Range('D44').Value = '1.443.31'
$ws.Range('E45').Value = '  -1.06%  '
$ws.Range('E46').Value = '  +6.91%  '
$ws.Range('E47').Value = '  +2.66%  '
$ws.Range('E48').Value = '  +0.82%  '
Set-TextValue 'D49' '7.41'
$ws.Range('E49').Value = '  +2.14%  '
Set-TextValue 'D50' '3.00'
$ws.Range('E50').Value = '  -0.97%  '
$ws.Range('D51').Value = '2.264.83'
$ws.Range('E51').Value = '  -1.46%  '
